$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# ColumnWidth (character units) stores internally with a fixed +5/6 pixel
# padding offset, so subtract that offset to land exactly on the target
# OOXML <col width="..."> value.
$offset = 0.8333333333333333
$ws.Columns("C").ColumnWidth = 63 - $offset
$ws.Columns("D").ColumnWidth = 33 - $offset
$ws.Columns("F").ColumnWidth = 17 - $offset
$ws.Columns("H").ColumnWidth = 45 - $offset

# Column A holds numeric-looking IDs that must stay text (inline/shared string),
# so force text format before writing any values into it.
$ws.Range("A2:A10").NumberFormat = "@"

# --- Row 2 (updated) ---
$ws.Range("A2").Value = "1328931"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328931"
$ws.Range("C2").Value = "Market Research of Wine and Spirits"
$ws.Range("D2").Value = "日本、大分県別府市"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "1 applicant"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "LINES Co., Ltd."

# --- Row 3 (updated) ---
$ws.Range("A3").Value = "1328849"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328849"
$ws.Range("C3").Value = "Geospatial Data Processing Intern"
$ws.Range("D3").Value = "Novi Sad, Serbia"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "1 applicant"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "DataDEV"

# --- Row 4 (updated) ---
$ws.Range("A4").Value = "1328630"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328630"
$ws.Range("C4").Value = "Marketing Intern"
$ws.Range("D4").Value = "Hyderabad, Telangana, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Amaavi Luxe Travels"

# --- Row 5 (updated) ---
$ws.Range("A5").Value = "1328615"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328615"
$ws.Range("C5").Value = "Back Office Planner"
$ws.Range("D5").Value = "Madrid, Spain"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "124 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "Mitsubishi Power Europe Sucursal en España"

# --- Row 6 (updated) ---
$ws.Range("A6").Value = "1328614"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1328614"
$ws.Range("C6").Value = "Field Service Engineer"
$ws.Range("D6").Value = "Madrid, Spain"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "64 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "Mitsubishi Power Europe Sucursal en España"

# --- Row 7 (new) ---
$ws.Range("A7").Value = "1328185"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1328185"
$ws.Range("C7").Value = "Data Scientist"
$ws.Range("D7").Value = "Frankfurt am Main, Deutschland"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "168 applicants"
$ws.Range("G7").Value = "3 - 6 Months"
$ws.Range("H7").Value = "Greyfood GmbH"

# --- Row 8 (new) ---
$ws.Range("A8").Value = "1327813"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327813"
$ws.Range("C8").Value = "Nursery Spanish Practitioner"
$ws.Range("D8").Value = "Ashby-de-la-Zouch LE65, UK"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "31 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "Bilingual Day Nursery and Preschool Ltd"

# --- Row 9 (new) ---
$ws.Range("A9").Value = "1327006"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1327006"
$ws.Range("C9").Value = "[Impact Florianópolis] Global HR Innovator Internship"
$ws.Range("D9").Value = "Balneário Camboriú, SC, Brasil"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "52 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "WTM do Brasil"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "1326913"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1326913"
$ws.Range("C10").Value = "Occupational Health and Safety Project Specialist (Mine Opp)"
$ws.Range("D10").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "28 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "Sodexo Mexico"
